# Corrected most names to the official names from website.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3-14: "Kalaburagi" -> "Kalaburagi (Gulbarga)" in column G
foreach ($r in 3..14) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "Kalaburagi") {
        $cell.Value2 = "Kalaburagi (Gulbarga)"
    }
}

# Rows 28-47: "Yadgiri" -> "Yadgir" in column G (rows 36 and 40 hold different text and are untouched)
foreach ($r in 28..47) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "Yadgiri") {
        $cell.Value2 = "Yadgir"
    }
}

# Remove the stray empty inline-string cells at F36 and F40
$ws.Cells.Item(36, 6).ClearContents()
$ws.Cells.Item(40, 6).ClearContents()
